$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.989.78'
$ws.Range('E2').Value = '''  +6.27%  '
$ws.Range('D3').Value = '''2.231.84'
$ws.Range('E3').Value = '''  +3.28%  '
$ws.Range('D5').Value = '''231.41'
$ws.Range('E5').Value = '''  +1.99%  '
$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '''  +0.64%  '
$ws.Range('D7').Value = '''60.89'
$ws.Range('E7').Value = '''  -2.70%  '
$ws.Range('E8').Value = '''  +0.01%  '
$ws.Range('E9').Value = '''  +3.31%  '
$ws.Range('D10').Value = '''58.95'
$ws.Range('E10').Value = '''  +1.11%  '
$ws.Range('D11').Value = '''0.0890'
$ws.Range('E11').Value = '''  +5.40%  '
$ws.Range('E12').Value = '''  -0.11%  '
$ws.Range('D13').Value = '''2.559.46'
$ws.Range('E13').Value = '''  +3.08%  '
$ws.Range('D14').Value = '''15.66'
$ws.Range('E14').Value = '''  -0.89%  '
$ws.Range('D15').Value = '''21.81'
$ws.Range('E15').Value = '''  +0.77%  '
$ws.Range('D16').Value = '''0.798'
$ws.Range('E16').Value = '''  -0.68%  '
$ws.Range('D17').Value = '''5.56'
$ws.Range('E17').Value = '''  +2.21%  '
$ws.Range('D18').Value = '''2.231.15'
$ws.Range('E18').Value = '''  +2.97%  '
$ws.Range('D19').Value = '''41.856.20'
$ws.Range('E19').Value = '''  +5.96%  '
$ws.Range('D20').Value = '''72.74'
$ws.Range('E20').Value = '''  +1.69%  '
$ws.Range('D21').Value = '''0.0₃0892'
$ws.Range('E21').Value = '''  +0.77%  '
$ws.Range('D22').Value = '''6.04'
$ws.Range('E22').Value = '''  +0.80%  '
$ws.Range('D23').Value = '''249.58'
$ws.Range('E23').Value = '''  +9.84%  '
$ws.Range('E24').Value = '''  -0.02%  '
$ws.Range('E25').Value = '''  +1.62%  '
$ws.Range('D26').Value = '''2.32'
$ws.Range('E26').Value = '''  -0.04%  '
$ws.Range('D27').Value = '''9.67'
$ws.Range('E27').Value = '''  +2.69%  '
$ws.Range('D28').Value = '''0.143'
$ws.Range('E28').Value = '''  +3.87%  '
$ws.Range('D29').Value = '''167.46'
$ws.Range('E29').Value = '''  -1.62%  '
$ws.Range('D30').Value = '''19.94'
$ws.Range('E30').Value = '''  +1.71%  '
$ws.Range('E31').Value = '''  -2.08%  '
$ws.Range('E32').Value = '''  -1.92%  '
$ws.Range('D33').Value = '''0.122'
$ws.Range('E33').Value = '''  +0.09%  '
$ws.Range('D34').Value = '''4.95'
$ws.Range('E34').Value = '''  +5.76%  '
$ws.Range('D35').Value = '''4.62'
$ws.Range('E35').Value = '''  +3.59%  '
$ws.Range('D36').Value = '''0.0628'
$ws.Range('E36').Value = '''  +2.05%  '
$ws.Range('E37').Value = '''  -4.07%  '
$ws.Range('D38').Value = '''3.69'
$ws.Range('E38').Value = '''  -2.90%  '
$ws.Range('D39').Value = '''2.36'
$ws.Range('E39').Value = '''  -1.06%  '
$ws.Range('E40').Value = '''  +31.12%  '
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '''  +0.00%  '
$ws.Range('D42').Value = '''0.0240'
$ws.Range('E42').Value = '''  +6.04%  '
$ws.Range('D43').Value = '''4.80'
$ws.Range('E43').Value = '''  -3.10%  '
$ws.Range('D44').Value = '''8.54'
$ws.Range('E44').Value = '''  +8.58%  '
$ws.Range('E45').Value = '''  +7.43%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''99.05'
$ws.Range('E46').Value = '''  -2.76%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '''1.21'
$ws.Range('E47').Value = '''  +0.90%  '
$ws.Range('D48').Value = '''1.470.92'
$ws.Range('E48').Value = '''  -2.64%  '
$ws.Range('D49').Value = '''16.53'
$ws.Range('E49').Value = '''  -6.53%  '
$ws.Range('E50').Value = '''  +0.16%  '
$ws.Range('D51').Value = '''1.08'
$ws.Range('E51').Value = '''  -0.80%  '
